# Update NATMI LR-pair sheet (Hcrt-Hcrtr1) with the new TPM-derived results.
# The "Resolving-Mac" cluster from the old run has been replaced by two
# clusters in the new run ("ECs" and "Inflammatory-Mac"), so the 3 original
# data rows (one per target cluster) are replaced by 8 rows: one per
# (sending cluster x target cluster) combination for the two sending
# clusters "ECs" and "Inflammatory-Mac", each against the four target
# clusters ECs / FAPs / Inflammatory-Mac / MuSCs.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @("ECs",              "Hcrt", "Hcrtr1", "ECs",              1, 0.3333333333333333, 0.1929206666666667, 0.578762,            0.9049887103533253, 0.9049887103533253, 2, 0.6666666666666666, 0.1991803333333333, 0.597541,            0.4922665899964288, 0.5074072124153813, 0.03842600269355555,  0.345834024242,   0.4454957064308973,   0.4591977987877717),
    @("ECs",              "Hcrt", "Hcrtr1", "FAPs",             1, 0.3333333333333333, 0.1929206666666667, 0.578762,            0.9049887103533253, 0.9049887103533253, 2, 0.6666666666666666, 0.1525796666666667, 0.457739,            0.3770948213400843, 0.3886931106046351, 0.02943577101311111,  0.264921939118,   0.3412665560454805,   0.3517628768893111),
    @("ECs",              "Hcrt", "Hcrtr1", "Inflammatory-Mac", 1, 0.3333333333333333, 0.1929206666666667, 0.578762,            0.9049887103533253, 0.9049887103533253, 1, 0.3333333333333333, 0.01663833333333333,0.049915,            0.04112100565429275, 0.04238576266350554, 0.003209878358888889, 0.02888890523,    0.03721404587551019,  0.038358636690188),
    @("ECs",              "Hcrt", "Hcrtr1", "MuSCs",            1, 0.3333333333333333, 0.1929206666666667, 0.578762,            0.9049887103533253, 0.9049887103533253, 1, 0.5,                0.0362205,           0.07244100000000001, 0.08951758300919427, 0.0615139143164781, 0.006987683007,       0.041926098042,   0.08101240200143746,  0.05566939798605447),
    @("Inflammatory-Mac", "Hcrt", "Hcrtr1", "ECs",              1, 0.3333333333333333, 0.020254,           0.060762,            0.09501128964667471,0.09501128964667471,2, 0.6666666666666666, 0.1991803333333333, 0.597541,            0.4922665899964288, 0.5074072124153813, 0.004034198471333333, 0.036307786242,   0.04677088356553156,  0.0482094136276096),
    @("Inflammatory-Mac", "Hcrt", "Hcrtr1", "FAPs",             1, 0.3333333333333333, 0.020254,           0.060762,            0.09501128964667471,0.09501128964667471,2, 0.6666666666666666, 0.1525796666666667, 0.457739,            0.3770948213400843, 0.3886931106046351, 0.003090348568666666, 0.027813137118,   0.0358282652946038,   0.03693023371532395),
    @("Inflammatory-Mac", "Hcrt", "Hcrtr1", "Inflammatory-Mac", 1, 0.3333333333333333, 0.020254,           0.060762,            0.09501128964667471,0.09501128964667471,1, 0.3333333333333333, 0.01663833333333333,0.049915,            0.04112100565429275, 0.04238576266350554, 0.0003369928033333333,0.00303293523,    0.003906959778782557, 0.004027125973317535),
    @("Inflammatory-Mac", "Hcrt", "Hcrtr1", "MuSCs",            1, 0.3333333333333333, 0.020254,           0.060762,            0.09501128964667471,0.09501128964667471,1, 0.5,                0.0362205,           0.07244100000000001, 0.08951758300919427, 0.0615139143164781, 0.000733610007,       0.004401660042,   0.008505181007756803, 0.005844516330423631)
)

$r = 2
foreach ($row in $rows) {
    $ws.Cells.Item($r, 1).Value  = $row[0]
    $ws.Cells.Item($r, 2).Value  = $row[1]
    $ws.Cells.Item($r, 3).Value  = $row[2]
    $ws.Cells.Item($r, 4).Value  = $row[3]
    $ws.Cells.Item($r, 5).Value  = $row[4]
    $ws.Cells.Item($r, 6).Value  = $row[5]
    $ws.Cells.Item($r, 7).Value  = $row[6]
    $ws.Cells.Item($r, 8).Value  = $row[7]
    $ws.Cells.Item($r, 9).Value  = $row[8]
    $ws.Cells.Item($r, 10).Value = $row[9]
    $ws.Cells.Item($r, 11).Value = $row[10]
    $ws.Cells.Item($r, 12).Value = $row[11]
    $ws.Cells.Item($r, 13).Value = $row[12]
    $ws.Cells.Item($r, 14).Value = $row[13]
    $ws.Cells.Item($r, 15).Value = $row[14]
    $ws.Cells.Item($r, 16).Value = $row[15]
    $ws.Cells.Item($r, 17).Value = $row[16]
    $ws.Cells.Item($r, 18).Value = $row[17]
    $ws.Cells.Item($r, 19).Value = $row[18]
    $ws.Cells.Item($r, 20).Value = $row[19]
    $r = $r + 1
}
